$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.599.08"
$ws.Range("E2").Value = "  -3.28%  "
$ws.Range("D3").Value = "1.850.42"
$ws.Range("E3").Value = "  -3.85%  "
$ws.Range("E4").Value = "  -1.03%  "
$ws.Range("D5").Value = "'334.24"
$ws.Range("E5").Value = "  +2.38%  "
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("D7").Value = "'0.4656"
$ws.Range("E7").Value = "  -3.47%  "
$ws.Range("D8").Value = "'0.3921"
$ws.Range("E8").Value = "  -3.74%  "
$ws.Range("D9").Value = "'46.52"
$ws.Range("E9").Value = "  -2.43%  "
$ws.Range("D10").Value = "'0.07907"
$ws.Range("E10").Value = "  -4.10%  "
$ws.Range("D11").Value = "'0.9860"
$ws.Range("E11").Value = "  -2.58%  "
$ws.Range("D12").Value = "'22.24"
$ws.Range("E12").Value = "  -5.90%  "
$ws.Range("D13").Value = "2.049.49"
$ws.Range("E13").Value = "  +7.63%  "
$ws.Range("D14").Value = "'5.854"
$ws.Range("E14").Value = "  -3.89%  "
$ws.Range("E15").Value = "  -3.60%  "
$ws.Range("D16").Value = "'0.06862"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "'87.88"
$ws.Range("E17").Value = "  -4.29%  "
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").Value = "'0.00001004"
$ws.Range("E19").Value = "  -3.43%  "
$ws.Range("D20").Value = "'17.14"
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").Value = "28.627.43"
$ws.Range("E22").Value = "  -3.15%  "
$ws.Range("D23").Value = "'5.408"
$ws.Range("E23").Value = "  -4.87%  "
$ws.Range("D24").Value = "'11.30"
$ws.Range("E24").Value = "  -5.20%  "
$ws.Range("D25").Value = "2.184.01"
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("D26").Value = "'2.126"
$ws.Range("E26").Value = "  -2.74%  "
$ws.Range("D27").Value = "'153.18"
$ws.Range("E27").Value = "  -1.85%  "
$ws.Range("D28").Value = "'19.44"
$ws.Range("E28").Value = "  -2.96%  "
$ws.Range("D29").Value = "'6.088"
$ws.Range("E29").Value = "  -6.08%  "
$ws.Range("D30").Value = "'2.016"
$ws.Range("E30").Value = "  -4.00%  "
$ws.Range("D31").Value = "'117.51"
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("D32").Value = "'0.9812"
$ws.Range("E32").Value = "  -3.63%  "
$ws.Range("D33").Value = "'0.09435"
$ws.Range("E33").Value = "  -2.22%  "
$ws.Range("D34").Value = "'5.369"
$ws.Range("D35").Value = "'3.482"
$ws.Range("E35").Value = "  -2.08%  "
$ws.Range("D36").Value = "'1.349"
$ws.Range("E36").Value = "  -2.14%  "
$ws.Range("D37").Value = "'0.06151"
$ws.Range("E37").Value = "  -3.55%  "
$ws.Range("D38").Value = "'0.02199"
$ws.Range("E38").Value = "  -4.30%  "
$ws.Range("E39").Value = "  -2.72%  "
$ws.Range("D40").Value = "'0.5710"
$ws.Range("E40").Value = "  -4.24%  "
$ws.Range("D41").Value = "'7.618"
$ws.Range("E41").Value = "  -3.59%  "
$ws.Range("D42").Value = "'10.14"
$ws.Range("E42").Value = "  -6.17%  "
$ws.Range("D43").Value = "'0.1799"
$ws.Range("E43").Value = "  -2.83%  "
$ws.Range("D44").Value = "'2.370"
$ws.Range("E44").Value = "  -4.12%  "
$ws.Range("D45").Value = "'1.249"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("D46").Value = "'11.87"
$ws.Range("E46").Value = "  -4.69%  "
$ws.Range("D47").Value = "'0.5396"
$ws.Range("E47").Value = "  -3.29%  "
$ws.Range("D48").Value = "'0.07162"
$ws.Range("E48").Value = "  -4.54%  "
$ws.Range("D49").Value = "'1.911"
$ws.Range("E49").Value = "  -2.24%  "
$ws.Range("D50").Value = "'114.14"
$ws.Range("E50").Value = "  -4.28%  "
$ws.Range("D51").Value = "'42.70"
$ws.Range("E51").Value = "  +1.73%  "
